# Hands On Demos - Day 1.
# Remove the stray "object 14" freeform shape from slide 5.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "object 14" -and $sh.Id -eq 14) {
        $sh.Delete()
        break
    }
}
